$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 456.7857
$ws.Range("I55").Value = 30.5
$ws.Range("J55").Value = 527.8333
$ws.Range("K55").Value = 30.5
$ws.Range("L55").Value = 527.8333
$ws.Range("M55").Value = 183.5
$ws.Range("N55").Value = -955.8333
$ws.Range("H112").Value = 4499.788
$ws.Range("J112").Value = 4773.968
$ws.Range("L112").Value = 14321.904
$ws.Range("N112").Value = -16537.904
$ws.Range("H116").Value = 13501
$ws.Range("J116").Value = 3399.8
$ws.Range("L116").Value = 3399.8
$ws.Range("N116").Value = -10283.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 1152.8
$ws.Range("I74").Value = 1127.9048
$ws.Range("J74").Value = 1210.8889
$ws.Range("K74").Value = 1127.9048
$ws.Range("L74").Value = 1210.8889
$ws.Range("M74").Value = -253.9048
$ws.Range("N74").Value = -2958.8889
$ws.Range("H77").Value = 1152.8
$ws.Range("I77").Value = 1127.9048
$ws.Range("J77").Value = 1210.8889
$ws.Range("K77").Value = 5639.524
$ws.Range("L77").Value = 6054.4445
$ws.Range("M77").Value = -1271.524
$ws.Range("N77").Value = -14790.4445
$ws.Range("H122").Value = 29823.834
$ws.Range("I122").Value = 1939.8928
$ws.Range("J122").Value = 127417.625
$ws.Range("K122").Value = 5819.678400000001
$ws.Range("L122").Value = 382252.875
$ws.Range("M122").Value = -3369.678400000001
$ws.Range("N122").Value = -387152.875
$ws.Range("H132").Value = 20430252
$ws.Range("I132").Value = 26317284
$ws.Range("J132").Value = 93231.63
$ws.Range("K132").Value = 78951852
$ws.Range("L132").Value = 279694.89
$ws.Range("M132").Value = -78949322
$ws.Range("N132").Value = -284754.89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 514.3333
$ws.Range("I80").Value = 891.5714
$ws.Range("J80").Value = 325.7143
$ws.Range("K80").Value = 891.5714
$ws.Range("L80").Value = 325.7143
$ws.Range("M80").Value = 106.4286
$ws.Range("N80").Value = -2321.7143
$ws.Range("H83").Value = 514.3333
$ws.Range("I83").Value = 891.5714
$ws.Range("J83").Value = 325.7143
$ws.Range("K83").Value = 4457.857
$ws.Range("L83").Value = 1628.5715
$ws.Range("M83").Value = 534.143
$ws.Range("N83").Value = -11612.5715
$ws.Range("H132").Value = 18000
$ws.Range("J132").Value = 18000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -28120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1212.4445
$ws.Range("I58").Value = 1287.4286
$ws.Range("J58").Value = 950
$ws.Range("K58").Value = 1287.4286
$ws.Range("L58").Value = 950
$ws.Range("M58").Value = -1084.4286
$ws.Range("N58").Value = -1356
$ws.Range("H94").Value = 6814.727
$ws.Range("I94").Value = 1282.4
$ws.Range("J94").Value = 11425
$ws.Range("K94").Value = 1282.4
$ws.Range("L94").Value = 11425
$ws.Range("M94").Value = -831.4000000000001
$ws.Range("N94").Value = -12327
$ws.Range("H122").Value = 836.96295
$ws.Range("I122").Value = 511.5
$ws.Range("J122").Value = 1310.3636
$ws.Range("K122").Value = 1534.5
$ws.Range("L122").Value = 3931.0908
$ws.Range("M122").Value = 915.5
$ws.Range("N122").Value = -8831.0908
$ws.Range("H136").Value = 1212.4445
$ws.Range("I136").Value = 1287.4286
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 3862.2858
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -1312.2858
$ws.Range("N136").Value = -7950
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 498.7
$ws.Range("I114").Value = 396.4
$ws.Range("J114").Value = 601
$ws.Range("K114").Value = 1189.2
$ws.Range("L114").Value = 1803
$ws.Range("M114").Value = 2064.8
$ws.Range("N114").Value = -8311
$ws.Range("H131").Value = 5377206.5
$ws.Range("I131").Value = 370.5
$ws.Range("J131").Value = 5883261.5
$ws.Range("K131").Value = 1111.5
$ws.Range("L131").Value = 17649784.5
$ws.Range("M131").Value = 3928.5
$ws.Range("N131").Value = -17659864.5
$ws.Range("H132").Value = 2299.9167
$ws.Range("I132").Value = 801.3333
$ws.Range("K132").Value = 7211.9997
$ws.Range("M132").Value = -4681.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 12.333333
$ws.Range("I2").Value = 12.8
$ws.Range("K2").Value = 12.8
$ws.Range("M2").Value = 100.2
$ws.Range("H21").Value = 2021200
$ws.Range("I21").Value = 5050000
$ws.Range("K21").Value = 5050000
$ws.Range("M21").Value = -5049827
$ws.Range("H30").Value = 2021200
$ws.Range("I30").Value = 5050000
$ws.Range("K30").Value = 5050000
$ws.Range("M30").Value = -5049895
$ws.Range("H122").Value = 1585.875
$ws.Range("I122").Value = 1397.8182
$ws.Range("J122").Value = 1999.6
$ws.Range("K122").Value = 4193.4546
$ws.Range("L122").Value = 5998.799999999999
$ws.Range("M122").Value = -1743.4546
$ws.Range("N122").Value = -10898.8
$ws.Range("H124").Value = 58672
$ws.Range("J124").Value = 58672
$ws.Range("L124").Value = 58672
$ws.Range("N124").Value = -68492
$ws.Range("H128").Value = 59769.668
$ws.Range("J128").Value = 59769.668
$ws.Range("L128").Value = 59769.668
$ws.Range("N128").Value = -69729.66800000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1889.5
$ws.Range("I40").Value = 1668.08
$ws.Range("K40").Value = 1668.08
$ws.Range("M40").Value = -1532.08
$ws.Range("H100").Value = 2793.3428
$ws.Range("I100").Value = 1271.8
$ws.Range("J100").Value = 3046.9333
$ws.Range("K100").Value = 1271.8
$ws.Range("L100").Value = 3046.9333
$ws.Range("M100").Value = -730.8
$ws.Range("N100").Value = -4128.933300000001
$ws.Range("H122").Value = 2352.8462
$ws.Range("J122").Value = 2489.3333
$ws.Range("L122").Value = 7467.999899999999
$ws.Range("N122").Value = -12367.9999
$ws.Range("H133").Value = 101041.664
$ws.Range("J133").Value = 101041.664
$ws.Range("L133").Value = 101041.664
$ws.Range("N133").Value = -106101.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1228.5
$ws.Range("I122").Value = 1166
$ws.Range("J122").Value = 1441
$ws.Range("K122").Value = 3498
$ws.Range("L122").Value = 4323
$ws.Range("M122").Value = -1048
$ws.Range("N122").Value = -9223
